$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D column (last-refresh timestamp) for rows 2-47 uniformly
$ws.Range("D2:D47").Value = 46019.384479166663

# Update rows 18-47 with refreshed station/terminal/charge-time data
$ws.Range("A18").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B18").Value = "101号直流"
$ws.Range("C18").Value = 46013.540347222224

$ws.Range("A19").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B19").Value = "705号直流"
$ws.Range("C19").Value = 46016.576261574075

$ws.Range("A20").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B20").Value = "001A号直流"
$ws.Range("C20").Value = 46017.034155092595

$ws.Range("A21").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B21").Value = "406号直流"
$ws.Range("C21").Value = 46017.04790509259

$ws.Range("A22").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B22").Value = "404号直流"
$ws.Range("C22").Value = 46017.677615740744

$ws.Range("A23").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B23").Value = "804号直流"
$ws.Range("C23").Value = 46017.83997685185

$ws.Range("A24").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B24").Value = "604号直流"
$ws.Range("C24").Value = 46017.885405092595

$ws.Range("A25").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B25").Value = "A01号直流"
$ws.Range("C25").Value = 46017.975324074076

$ws.Range("A26").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B26").Value = "107号直流"
$ws.Range("C26").Value = 46018.464791666665

$ws.Range("A27").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B27").Value = "305号直流"
$ws.Range("C27").Value = 46018.5165162037

$ws.Range("A28").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B28").Value = "805号直流"
$ws.Range("C28").Value = 46018.52658564815

$ws.Range("A29").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B29").Value = "003B号直流"
$ws.Range("C29").Value = 46018.535995370374

$ws.Range("A30").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B30").Value = "902号直流"
$ws.Range("C30").Value = 46018.54603009259

$ws.Range("A31").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B31").Value = "402号直流"
$ws.Range("C31").Value = 46018.55363425926

$ws.Range("A32").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B32").Value = "206号直流"
$ws.Range("C32").Value = 46018.55541666667

$ws.Range("A33").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B33").Value = "203号直流"
$ws.Range("C33").Value = 46018.56657407407

$ws.Range("A34").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B34").Value = "503号直流"
$ws.Range("C34").Value = 46018.58

$ws.Range("A35").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B35").Value = "703号直流"
$ws.Range("C35").Value = 46018.58696759259

$ws.Range("A36").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B36").Value = "501号直流"
$ws.Range("C36").Value = 46018.595717592594

$ws.Range("A37").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B37").Value = "B01号直流"
$ws.Range("C37").Value = 46018.60731481481

$ws.Range("A38").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B38").Value = "101号直流"
$ws.Range("C38").Value = 46018.62708333333

$ws.Range("A39").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B39").Value = "105号直流"
$ws.Range("C39").Value = 46018.639502314814

$ws.Range("A40").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B40").Value = "002A号直流"
$ws.Range("C40").Value = 46018.64497685185

$ws.Range("A41").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B41").Value = "903号直流"
$ws.Range("C41").Value = 46018.65012731482

$ws.Range("A42").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B42").Value = "905号直流"
$ws.Range("C42").Value = 46018.66664351852

$ws.Range("A43").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B43").Value = "401号直流"
$ws.Range("C43").Value = 46018.67826388889

$ws.Range("A44").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B44").Value = "505号直流"
$ws.Range("C44").Value = 46018.70648148148

$ws.Range("A45").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B45").Value = "105号直流"
$ws.Range("C45").Value = 46018.72219907407

$ws.Range("A46").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B46").Value = "802号直流"
$ws.Range("C46").Value = 46018.73488425926

$ws.Range("A47").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B47").Value = "603号直流"
$ws.Range("C47").Value = 46018.77072916667

# Clear rows 48-64 (no longer have data, but styles remain)
$ws.Range("A48:D64").ClearContents()

# Remove row 65 entirely (dataset shrank by one entry)
$ws.Range("A65").EntireRow.Delete()

# Restore active selection to reflect the last user interaction
$ws.Range("G15").Select()
